# Applies the Wed Oct 11 07:36:11 UTC 2023 cryptos-list refresh described in the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.108.37"
$ws.Range("E2").Value = "  -2.29%  "
$ws.Range("D3").Value = "1.559.83"
$ws.Range("E3").Value = "  -2.32%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "206.60"
$ws.Range("E5").Value = "  -0.96%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.487"
$ws.Range("E6").Value = "  -3.11%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.11"
$ws.Range("E8").Value = "  -1.28%  "
$ws.Range("E9").Value = "  -2.51%  "
$ws.Range("E10").Value = "  -0.56%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0865"
$ws.Range("E11").Value = "  -0.47%  "
$ws.Range("D12").Value = "1.780.97"
$ws.Range("E12").Value = "  -2.29%  "
$ws.Range("D13").Value = "1.562.40"
$ws.Range("E13").Value = "  -2.10%  "
$ws.Range("E14").Value = "  -2.86%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.516"
$ws.Range("E15").Value = "  -3.50%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "62.88"
$ws.Range("E16").Value = "  -1.14%  "
$ws.Range("D17").Value = "27.121.61"
$ws.Range("E17").Value = "  -2.20%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "214.95"
$ws.Range("E18").Value = "  -2.23%  "
$ws.Range("E19").Value = "  -1.86%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.22"
$ws.Range("E20").Value = "  -2.07%  "
$ws.Range("E21").Value = "  -0.01%  "
$ws.Range("E22").Value = "  -1.13%  "
$ws.Range("E23").Value = "  -4.68%  "
$ws.Range("E24").Value = "  +0.19%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.67"
$ws.Range("E25").Value = "  -1.73%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.59"
$ws.Range("E26").Value = "  -2.78%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "14.90"
$ws.Range("E27").Value = "  -1.78%  "
$ws.Range("E28").Value = "  +0.03%  "
$ws.Range("E29").Value = "  -1.68%  "
$ws.Range("E30").Value = "  -1.71%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0461"
$ws.Range("E31").Value = "  -2.47%  "
$ws.Range("E32").Value = "  -2.49%  "
$ws.Range("D33").Value = "1.383.25"
$ws.Range("E33").Value = "  +0.31%  "
$ws.Range("E34").Value = "  -1.55%  "
$ws.Range("E35").Value = "  -0.66%  "
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.28"
$ws.Range("E36").Value = "  -1.66%  "
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.943"
$ws.Range("E37").Value = "  -2.96%  "
$ws.Range("E38").Value = "  -1.95%  "
$ws.Range("E39").Value = "  -2.50%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.513"
$ws.Range("E40").Value = "  -4.43%  "
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("E42").Value = "  +1.36%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.78"
$ws.Range("E43").Value = "  +2.65%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "63.32"
$ws.Range("E44").Value = "  -2.11%  "
$ws.Range("E45").Value = "  -0.25%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.25"
$ws.Range("E46").Value = "  +0.50%  "
$ws.Range("D47").Value = "1.694.19"
$ws.Range("E47").Value = "  -2.22%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "85.40"
$ws.Range("E48").Value = "  -1.87%  "
$ws.Range("D49").Value = "0.0₇0984"
$ws.Range("E49").Value = "  -3.18%  "
$ws.Range("E50").Value = "  -0.98%  "
